# Revert "Powerpoint writer: consolidate text run nodes."
# Split runs that end in a trailing space back into a word-run + a
# separate space-run, matching the pre-consolidation OOXML shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1: "Testing " + "custom " + "properties"
#     -> "Testing" / " " / "custom" / " " / "properties"
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$word1 = $titleRange.Characters(1, 7)
$word1.Text = "Testing"

$word2 = $titleRange.Characters(9, 6)
$word2.Text = "custom"

# --- Subtitle 2: "A. " + "M." -> "A." / " " / "M."
$subShape = $s.Shapes.Item(2)
$subRange = $subShape.TextFrame.TextRange

$word3 = $subRange.Characters(3, 2)
$word3.Text = "A."
